{"js": "// Apply the four r\u00e9sum\u00e9 text updates described by the diff:\n//  1. Phone number line gets a \"Tel\u00e9fono: \" prefix.\n//  2. Current role line gets employer prefix + retitled role + \"Presente\".\n//  3. Previous role line gets employer prefix + retitled role + abbreviated dates.\n//  4. Earliest role line gets employer prefix + retitled role + abbreviated dates.\n\nconst replacements = [\n  {\n    find: \"(123) 456-7890\",\n    replacement: \"Tel\u00e9fono: (123) 456-7890\",\n  },\n  {\n    find: \"Animadora jefe (enero de 2018 - actualidad)\",\n    replacement: \"ABC Studios: Animator principal (enero de 2018 - Presente)\",\n  },\n  {\n    find: \"Animadora principal (junio de 2015 - diciembre de 2017)\",\n    replacement: \"XYZ Media: Animator Senior (jun 2015 - dic 2017)\",\n  },\n  {\n    find: \"Animadora j\u00fanior (septiembre de 2012 - mayo de 2015)\",\n    replacement: \"MNO Entertainment: Junior Animator (Sep 2012 - Mayo de 2015)\",\n  },\n];\n\nfor (const { find, replacement } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply the four r\u00e9sum\u00e9 text updates described by the diff:\n#  1. Phone number line gets a \"Tel\u00e9fono: \" prefix.\n#  2. Current role line gets employer prefix + retitled role + \"Presente\".\n#  3. Previous role line gets employer prefix + retitled role + abbreviated dates.\n#  4. Earliest role line gets employer prefix + retitled role + abbreviated dates.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $FindText\n    $find.Replacement.Text = $ReplaceText\n\n    $find.Execute(\n        $FindText,    # FindText\n        $false,       # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $ReplaceText, # ReplaceWith\n        2             # Replace (wdReplaceAll)\n    )\n}\n\nReplace-ExactText \"(123) 456-7890\" \"Tel\u00e9fono: (123) 456-7890\"\nReplace-ExactText \"Animadora jefe (enero de 2018 - actualidad)\" \"ABC Studios: Animator principal (enero de 2018 - Presente)\"\nReplace-ExactText \"Animadora principal (junio de 2015 - diciembre de 2017)\" \"XYZ Media: Animator Senior (jun 2015 - dic 2017)\"\nReplace-ExactText \"Animadora j\u00fanior (septiembre de 2012 - mayo de 2015)\" \"MNO Entertainment: Junior Animator (Sep 2012 - Mayo de 2015)\"\n"}
